$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 15462.667
$ws.Range("I62").Value = 11144
$ws.Range("J62").Value = 24100
$ws.Range("K62").Value = 11144
$ws.Range("L62").Value = 24100
$ws.Range("M62").Value = -10520
$ws.Range("N62").Value = -25348

# Row 65
$ws.Range("H65").Value = 15462.667
$ws.Range("I65").Value = 11144
$ws.Range("J65").Value = 24100
$ws.Range("K65").Value = 55720
$ws.Range("L65").Value = 120500
$ws.Range("M65").Value = -52600
$ws.Range("N65").Value = -126740

# Row 80
$ws.Range("H80").Value = 1144.7667
$ws.Range("I80").Value = 457
$ws.Range("J80").Value = 1394.8636
$ws.Range("K80").Value = 1371
$ws.Range("L80").Value = 4184.5908
$ws.Range("M80").Value = -373
$ws.Range("N80").Value = -6180.5908

# Row 83
$ws.Range("H83").Value = 1144.7667
$ws.Range("I83").Value = 457
$ws.Range("J83").Value = 1394.8636
$ws.Range("K83").Value = 4113
$ws.Range("L83").Value = 12553.7724
$ws.Range("M83").Value = 879
$ws.Range("N83").Value = -22537.7724

# Row 113
$ws.Range("H113").Value = 127375.5
$ws.Range("I113").Value = 252026.25
$ws.Range("J113").Value = 2724.75
$ws.Range("K113").Value = 252026.25
$ws.Range("L113").Value = 2724.75
$ws.Range("M113").Value = -248772.25
$ws.Range("N113").Value = -9232.75

# Row 132
$ws.Range("H132").Value = 285128.44
$ws.Range("I132").Value = 329385.78
$ws.Range("J132").Value = 80438.25
$ws.Range("K132").Value = 988157.3400000001
$ws.Range("L132").Value = 241314.75
$ws.Range("M132").Value = -985627.3400000001
$ws.Range("N132").Value = -246374.75

# Row 137
$ws.Range("H137").Value = 26317032
$ws.Range("I137").Value = 32259012
$ws.Range("K137").Value = 96777036
$ws.Range("M137").Value = -96774486

# Row 138
$ws.Range("H138").Value = 4264211.5
$ws.Range("I138").Value = 1192811.6
$ws.Range("J138").Value = 6292495
$ws.Range("K138").Value = 3578434.8
$ws.Range("L138").Value = 18877485
$ws.Range("M138").Value = -3573294.8
$ws.Range("N138").Value = -18887765

# Row 141
$ws.Range("H141").Value = 2645.587
$ws.Range("I141").Value = 1637.4706
$ws.Range("J141").Value = 5501.9165
$ws.Range("K141").Value = 4912.4118
$ws.Range("L141").Value = 16505.7495
$ws.Range("M141").Value = 267.5882000000001
$ws.Range("N141").Value = -26865.7495

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 15532.959
$ws.Range("I32").Value = 1962.6812
$ws.Range("K32").Value = 1962.6812
$ws.Range("M32").Value = -1675.6812

# Row 43
$ws.Range("H43").Value = 5850
$ws.Range("J43").Value = 5850
$ws.Range("L43").Value = 5850
$ws.Range("N43").Value = -6476

# Row 61
$ws.Range("H61").Value = 1939.7084
$ws.Range("I61").Value = 1260.5714
$ws.Range("J61").Value = 3768.1538
$ws.Range("K61").Value = 1260.5714
$ws.Range("L61").Value = 3768.1538
$ws.Range("M61").Value = -1048.5714
$ws.Range("N61").Value = -4192.1538

# Row 74
$ws.Range("H74").Value = 3211.1897
$ws.Range("I74").Value = 886.83673
$ws.Range("J74").Value = 15866
$ws.Range("K74").Value = 886.83673
$ws.Range("L74").Value = 15866
$ws.Range("M74").Value = -12.83672999999999
$ws.Range("N74").Value = -17614

# Row 77
$ws.Range("H77").Value = 3211.1897
$ws.Range("I77").Value = 886.83673
$ws.Range("J77").Value = 15866
$ws.Range("K77").Value = 4434.18365
$ws.Range("L77").Value = 79330
$ws.Range("M77").Value = -66.18364999999994
$ws.Range("N77").Value = -88066

# Row 103
$ws.Range("H103").Value = 500181
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 500181
$ws.Range("K103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("M103").Value = 500181
$ws.Range("N103").Value = -502525

# Row 136
$ws.Range("H136").Value = 1939.7084
$ws.Range("I136").Value = 1260.5714
$ws.Range("J136").Value = 3768.1538
$ws.Range("K136").Value = 3781.7142
$ws.Range("L136").Value = 11304.4614
$ws.Range("M136").Value = -1231.7142
$ws.Range("N136").Value = -16404.4614

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 8056.875
$ws.Range("I86").Value = 2237.875
$ws.Range("J86").Value = 13875.875
$ws.Range("K86").Value = 2237.875
$ws.Range("L86").Value = 13875.875
$ws.Range("M86").Value = -1114.875
$ws.Range("N86").Value = -16121.875

# Row 89
$ws.Range("H89").Value = 8056.875
$ws.Range("I89").Value = 2237.875
$ws.Range("J89").Value = 13875.875
$ws.Range("K89").Value = 11189.375
$ws.Range("L89").Value = 69379.375
$ws.Range("M89").Value = -5573.375
$ws.Range("N89").Value = -80611.375

# Row 102
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").ClearContents()
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = 0

# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# Row 134
$ws.Range("H134").Value = 2126.2144
$ws.Range("I134").Value = 1426.3695
$ws.Range("J134").Value = 3467.5833
$ws.Range("K134").Value = 4279.1085
$ws.Range("L134").Value = 10402.7499
$ws.Range("M134").Value = -1744.1085
$ws.Range("N134").Value = -15472.7499

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1542.2059
$ws.Range("I31").Value = 882.62964
$ws.Range("K31").Value = 882.62964
$ws.Range("M31").Value = -587.62964

# Row 34
$ws.Range("H34").Value = 1542.2059
$ws.Range("I34").Value = 882.62964
$ws.Range("K34").Value = 882.62964
$ws.Range("M34").Value = -680.62964

# Row 103
$ws.Range("H103").Value = 30000
$ws.Range("I103").Value = 30000
$ws.Range("K103").Value = 30000

# Row 122
$ws.Range("H122").Value = 2471.1538
$ws.Range("I122").Value = 1260.5714
$ws.Range("J122").Value = 3883.5
$ws.Range("K122").Value = 3781.7142
$ws.Range("L122").Value = 11650.5
$ws.Range("M122").Value = -1331.7142
$ws.Range("N122").Value = -16550.5

$ws = $wb.Worksheets.Item("CUL")
# Row 82
$ws.Range("H82").Value = 5000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 5000
$ws.Range("K82").Value = 0
$ws.Range("L82").ClearContents()
$ws.Range("M82").Value = 15000
$ws.Range("N82").Value = -15812

# Row 85
$ws.Range("H85").Value = 5000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 5000
$ws.Range("K85").Value = 0
$ws.Range("L85").ClearContents()
$ws.Range("M85").Value = 15000
$ws.Range("N85").Value = -17808

# Row 99
$ws.Range("H99").Value = 562.5
$ws.Range("I99").Value = 562.5
$ws.Range("K99").Value = 1687.5
$ws.Range("M99").Value = 558.5

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2803.1538
$ws.Range("I80").Value = 2644
$ws.Range("J80").Value = 3333.6667
$ws.Range("K80").Value = 2644
$ws.Range("L80").Value = 3333.6667
$ws.Range("M80").Value = -1646
$ws.Range("N80").Value = -5329.6667

# Row 83
$ws.Range("H83").Value = 2803.1538
$ws.Range("I83").Value = 2644
$ws.Range("J83").Value = 3333.6667
$ws.Range("K83").Value = 13220
$ws.Range("L83").Value = 16668.3335
$ws.Range("M83").Value = -8228
$ws.Range("N83").Value = -26652.3335

# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("N134").Value = 0

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 697.5714
$ws.Range("J55").Value = 716.8
$ws.Range("L55").Value = 716.8
$ws.Range("N55").Value = -1062.8

# Row 61
$ws.Range("H61").Value = 9847.315000000001
$ws.Range("I61").Value = 9443.6875
$ws.Range("K61").Value = 9443.6875
$ws.Range("M61").Value = -9241.6875

# Row 113
$ws.Range("H113").Value = 9847.315000000001
$ws.Range("I113").Value = 9443.6875
$ws.Range("K113").Value = 9443.6875
$ws.Range("M113").Value = -7273.6875

# Row 132
$ws.Range("H132").Value = 3330.7144
$ws.Range("I132").Value = 2215.4412
$ws.Range("J132").Value = 5858.6665
$ws.Range("K132").Value = 6646.323600000001
$ws.Range("L132").Value = 17575.9995
$ws.Range("M132").Value = -4116.323600000001
$ws.Range("N132").Value = -22635.9995

# Row 136
$ws.Range("H136").Value = 4027
$ws.Range("I136").Value = 2452.2104
$ws.Range("J136").Value = 12575.857
$ws.Range("K136").Value = 7356.6312
$ws.Range("L136").Value = 37727.571
$ws.Range("M136").Value = -4806.6312
$ws.Range("N136").Value = -42827.571

$ws = $wb.Worksheets.Item("WVR")
# Row 26
$ws.Range("H26").Value = 8402.4
$ws.Range("I26").Value = 6006
$ws.Range("J26").Value = 10000
$ws.Range("K26").Value = 6006
$ws.Range("L26").Value = 10000
$ws.Range("N26").Value = -10586

# Row 29
$ws.Range("H29").Value = 6527.5
$ws.Range("I29").Value = 1555
$ws.Range("J29").Value = 11500
$ws.Range("K29").Value = 1555
$ws.Range("L29").Value = 11500
$ws.Range("N29").Value = -12080

# Row 126
$ws.Range("H126").Value = 56647.277
$ws.Range("I126").Value = 111550.11
$ws.Range("J126").Value = 1744.4445
$ws.Range("K126").Value = 334650.33
$ws.Range("L126").Value = 5233.333500000001
$ws.Range("M126").Value = -332180.33
$ws.Range("N126").Value = -10173.3335

# Row 132
$ws.Range("H132").Value = 8623231
$ws.Range("I132").Value = 12502741
$ws.Range("J132").Value = 2098.2222
$ws.Range("K132").Value = 37508223
$ws.Range("L132").Value = 6294.6666
$ws.Range("M132").Value = -37505693
$ws.Range("N132").Value = -11354.6666

# Row 136
$ws.Range("H136").Value = 9834362
$ws.Range("I136").Value = 11145368
$ws.Range("J136").Value = 1822.5
$ws.Range("K136").Value = 1822.5
$ws.Range("L136").Value = 5467.5
$ws.Range("M136").Value = -33433554
$ws.Range("N136").Value = -10567.5
